$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Convert the "{m:if 1/0 = 42}" field into plain literal text.
# ------------------------------------------------------------------
$ifField = $d.Fields.Item(1)
$insertPos = $ifField.Code.Start - 1
$ifField.Delete()

$ins = $d.Range($insertPos, $insertPos)
$ifText = "{m:if 1/0 = 42}"
$ins.InsertBefore($ifText)

$newIfRange = $d.Range($insertPos, $insertPos + $ifText.Length)
$newIfRange.Font.Bold = 0
$newIfRange.Font.Color = -16777216

# ------------------------------------------------------------------
# 2) Update the stack trace text (prefix + line-number refresh).
# ------------------------------------------------------------------
$d.Content.Find.Execute("divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:", $true, $false, $false, $false, $false, $true, 1, $false, "    <---divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:", 2)

$d.Content.Find.Execute("at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1697)", $true, $false, $false, $false, $false, $true, 1, $false, "at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:2260)", 2)

$d.Content.Find.Execute("caseConditional(M2DocEvaluator.java:1)" + [char]10 + "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1459)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1684)", $true, $false, $false, $false, $false, $true, 1, $false, "caseConditional(M2DocEvaluator.java:1)" + [char]10 + "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2247)", 2)

$d.Content.Find.Execute("caseBlock(M2DocEvaluator.java:1)" + [char]10 + "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1459)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)", $true, $false, $false, $false, $false, $true, 1, $false, "caseBlock(M2DocEvaluator.java:1)" + [char]10 + "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:340)", 2)

$d.Content.Find.Execute("caseDocumentTemplate(M2DocEvaluator.java:1)" + [char]10 + "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1459)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)" + [char]10 + "	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)" + [char]10 + "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:506)" + [char]10 + "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:400)", $true, $false, $false, $false, $false, $true, 1, $false, "caseDocumentTemplate(M2DocEvaluator.java:1)" + [char]10 + "	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + [char]10 + "	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2022)" + [char]10 + "	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:324)" + [char]10 + "	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:912)" + [char]10 + "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:497)" + [char]10 + "	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)", 2)

$d.Content.Find.Execute("DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + "	at java.base/java.lang.reflect.Method.invoke(Method.java:568)", $true, $false, $false, $false, $false, $true, 1, $false, "DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + "	at java.base/java.lang.reflect.Method.invoke(Method.java:569)", 2)

$d.Content.Find.Execute("RemoteTestRunner.runTests(RemoteTestRunner.java:756)", $true, $false, $false, $false, $false, $true, 1, $false, "RemoteTestRunner.runTests(RemoteTestRunner.java:757)", 2)

# ------------------------------------------------------------------
# 3) Convert the two "{m:endif}" fields into plain literal text.
#    (Only the two "m:endif" fields remain at this point.)
# ------------------------------------------------------------------
$endField2 = $d.Fields.Item(2)
$pos2 = $endField2.Code.Start - 1
$endField2.Delete()
$ins2 = $d.Range($pos2, $pos2)
$endText = "{m:endif}"
$ins2.InsertBefore($endText)
$newEnd2 = $d.Range($pos2, $pos2 + $endText.Length)
$newEnd2.Font.Bold = 0
$newEnd2.Font.Color = -16777216

$endField1 = $d.Fields.Item(1)
$pos1 = $endField1.Code.Start - 1
$endField1.Delete()
$ins1 = $d.Range($pos1, $pos1)
$ins1.InsertBefore($endText)
$newEnd1 = $d.Range($pos1, $pos1 + $endText.Length)
$newEnd1.Font.Bold = 0
$newEnd1.Font.Color = -16777216
